# Updates symbol list values (Price and Volume(1h)) to reflect the latest
# market snapshot. Values must remain plain text (matching the sheet's
# existing inline-string cells) rather than being auto-converted to
# numbers/percentages by Excel, so we force a Text number format before
# assigning, then restore the default "Normal" style so no stray style
# index is left behind on the cell.
function Set-CellText {
    param($ws, $ref, $val)
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-CellText $ws "D2" "303.50"
Set-CellText $ws "E2" "1.25%"
Set-CellText $ws "D3" "32.17"
Set-CellText $ws "E3" "2.95%"
Set-CellText $ws "D4" "4.923"
Set-CellText $ws "E4" "-3.17%"
Set-CellText $ws "D5" "0.07867"
Set-CellText $ws "E5" "-1.10%"
Set-CellText $ws "D6" "2.048"
Set-CellText $ws "E6" "-6.87%"
Set-CellText $ws "D7" "7.846"
Set-CellText $ws "E7" "1.27%"
Set-CellText $ws "D8" "3.846"
Set-CellText $ws "E8" "-0.41%"
Set-CellText $ws "D9" "0.9244"
Set-CellText $ws "E9" "0.87%"
Set-CellText $ws "E10" "1.53%"
Set-CellText $ws "D11" "0.07820"
Set-CellText $ws "E11" "5.77%"
Set-CellText $ws "D12" "0.08602"
Set-CellText $ws "E12" "-7.98%"
Set-CellText $ws "D13" "0.03165"
Set-CellText $ws "E13" "4.85%"
Set-CellText $ws "E14" "0.34%"
Set-CellText $ws "D15" "0.001510"
Set-CellText $ws "E15" "-0.26%"
Set-CellText $ws "D16" "0.005796"
Set-CellText $ws "E16" "-1.98%"
Set-CellText $ws "E17" "2,110.07%"
Set-CellText $ws "D18" "3.465"
Set-CellText $ws "E18" "-0.36%"
Set-CellText $ws "D19" "2.116"
Set-CellText $ws "E19" "-6.58%"
Set-CellText $ws "D20" "0.3277"
Set-CellText $ws "E20" "0.16%"
Set-CellText $ws "D21" "0.1318"
Set-CellText $ws "E21" "-1.33%"
Set-CellText $ws "D22" "4.279"
Set-CellText $ws "E22" "2.97%"
Set-CellText $ws "D23" "0.1990"
Set-CellText $ws "E23" "17.21%"
Set-CellText $ws "D24" "0.04584"
Set-CellText $ws "E24" "-0.75%"
Set-CellText $ws "D25" "0.001224"
Set-CellText $ws "E25" "-1.46%"
Set-CellText $ws "D26" "0.004458"
Set-CellText $ws "E26" "-0.03%"
Set-CellText $ws "D27" "0.0001249"
Set-CellText $ws "E27" "4.22%"
Set-CellText $ws "E39" "-1.01%"
Set-CellText $ws "E40" "4.31%"
Set-CellText $ws "D41" "0.007474"
Set-CellText $ws "E41" "7.46%"
Set-CellText $ws "E42" "0.61%"
Set-CellText $ws "D43" "0.002358"
Set-CellText $ws "E43" "7.81%"
Set-CellText $ws "D44" "0.01044"
Set-CellText $ws "E44" "9.31%"
Set-CellText $ws "D45" "0.00006166"
Set-CellText $ws "E45" "-1.86%"
Set-CellText $ws "D46" "0.00000000749"
Set-CellText $ws "E46" "0.05%"
Set-CellText $ws "E47" "-61.10%"
Set-CellText $ws "D48" "0.8205"
Set-CellText $ws "E48" "9.85%"
Set-CellText $ws "D49" "0.00002099"
Set-CellText $ws "E49" "0.05%"
Set-CellText $ws "D50" "0.0001999"
Set-CellText $ws "E50" "0.05%"
